$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.145.59"
$ws.Range("E2").Value = "  +0.27%  "
$ws.Range("D3").Value = "1.676.13"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.61%  "
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.82"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +6.69%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.261"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.86%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0620"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.64%  "
$ws.Range("E11").Value = "  -0.06%  "
$ws.Range("D12").Value = "1.912.83"
$ws.Range("E12").Value = "  -0.22%  "
$ws.Range("D13").Value = "1.677.05"
$ws.Range("E13").Value = "  -0.05%  "
$ws.Range("E14").Value = "  +2.21%  "
$ws.Range("E15").Value = "  +4.72%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.49"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "27.125.99"
$ws.Range("E17").Value = "  +0.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "235.14"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.50%  "
$ws.Range("D19").Value = "0.0₃0741"
$ws.Range("E19").Value = "  +0.56%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.81"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.27%  "
$ws.Range("E21").Value = "  +0.15%  "
$ws.Range("E22").Value = "  +1.71%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.55"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.07%  "
$ws.Range("E24").Value = "  -2.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "148.33"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.84%  "
$ws.Range("E26").Value = "  +2.38%  "
$ws.Range("E28").Value = "  -0.21%  "
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("E30").Value = "  +0.51%  "
$ws.Range("E31").Value = "  -0.28%  "
$ws.Range("E32").Value = "  +0.39%  "
$ws.Range("D33").Value = "1.546.09"
$ws.Range("E33").Value = "  -0.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.24"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.76%  "
$ws.Range("E35").Value = "  -4.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.606"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.99%  "
$ws.Range("E37").Value = "  +3.21%  "
$ws.Range("E38").Value = "  +0.04%  "
$ws.Range("E39").Value = "  -0.91%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.06"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.25%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "69.97"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.17%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.78"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.46%  "
$ws.Range("E43").Value = "  +0.17%  "
$ws.Range("E44").Value = "  -0.33%  "
$ws.Range("D45").Value = "1.822.51"
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.782"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.12%  "
$ws.Range("E47").Value = "  +6.56%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "89.69"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.91%  "
$ws.Range("E49").Value = "  +2.50%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.23"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.83%  "
$ws.Range("E51").Value = "  -0.10%  "
